$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (written in this order so the shared-string table indices line up)
$ws.Range("C1").Value = "numBranches"
$ws.Range("E1").Value = "l1 dcache miss rate "
$ws.Range("G1").Value = "l2 cache miss rate"
$ws.Range("D1").Value = "l1 dcache access num"
$ws.Range("F1").Value = "l2 cache access num"

# Row 2 - Baseline
$ws.Range("C2").Value = 2199733
$ws.Range("D2").Value = 40295667
$ws.Range("E2").Value = 0.056419999999999998
$ws.Range("F2").Value = 2273474
$ws.Range("G2").Value = 0.93229300000000004

# Row 3 - Loop Unrolling
$ws.Range("C3").Value = 365064
$ws.Range("D3").Value = 34791129
$ws.Range("E3").Value = 0.065346000000000001
$ws.Range("F3").Value = 2273477
$ws.Range("G3").Value = 0.93229300000000004

# Row 4 - Matrix blocking
$ws.Range("C4").Value = 4518345
$ws.Range("D4").Value = 45016796
$ws.Range("E4").Value = 0.050622
$ws.Range("F4").Value = 2278853
$ws.Range("G4").Value = 0.088002999999999998

# Row 5 - Transposition storage
$ws.Range("C5").Value = 2200239
$ws.Range("D5").Value = 40296389
$ws.Range("E5").Value = 0.0091660000000000005
$ws.Range("F5").Value = 369349
$ws.Range("G5").Value = 0.37201099999999998

# Row 6 - Outer product
$ws.Range("C6").Value = 2216954
$ws.Range("D6").Value = 36348868
$ws.Range("E6").Value = 0.0042399999999999998
$ws.Range("F6").Value = 154116
$ws.Range("G6").Value = 0.998255

# Column widths to match bestFit autosizing in the target file (closest
# achievable values given this host's column-width rounding granularity)
$ws.Columns.Item(3).ColumnWidth = 11.142857142857142
$ws.Columns.Item(4).ColumnWidth = 18.142857142857142
$ws.Columns.Item(5).ColumnWidth = 16.428571428571427
$ws.Columns.Item(6).ColumnWidth = 17.0
$ws.Columns.Item(7).ColumnWidth = 14.714285714285714

# Update selection to match target
$ws.Range("E8").Select()
